$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values (columns B, C, D) for rows 2-9 ---
# Row 2 - ARDUINO
$ws.Range("B2").Value = 0.055
$ws.Range("C2").Value = 0.00036
$ws.Range("D2").Value = 5

# Row 3 - SIGFOX
$ws.Range("B3").Value = 0.049
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 3.3

# Row 4 - RTC (values unchanged, style still applied below)
$ws.Range("B4").Value = 0.0015
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 5

# Row 5 - GPS (values unchanged)
$ws.Range("B5").Value = 0.04
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 3.3

# Row 6 - DS18B20 (values unchanged)
$ws.Range("B6").Value = 0.0015
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 5

# Row 7 - DHT11 (values unchanged)
$ws.Range("B7").Value = 0.0005
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 5

# Row 8 - LCD
$ws.Range("B8").Value = 0.164
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 5

# Row 9 - HX711
$ws.Range("B9").Value = 0.0015
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 5

# --- Apply the "Satisfaisant" (green) cell style to all data rows 2-9 ---
# Columns A and D keep the default (General) number format under that style,
# columns B and C get the scientific number format (numFmtId 11 -> 0.00E+00)
$ws.Range("A2:A9").Style = "Satisfaisant"
$ws.Range("D2:D9").Style = "Satisfaisant"
$ws.Range("B2:C9").Style = "Satisfaisant"
$ws.Range("B2:C9").NumberFormat = "0.00E+00"

# --- Update the active cell selection ---
$ws.Range("G7").Select()
